$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")
$ws.Range("A10:G13").Insert()
Write-Host "inserted"
